$wb = $excel.ActiveWorkbook

$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Huoerxinhe Coal Mine, China, M1140, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)`""

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$newVersionString = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

for ($row = 2; $row -le 8; $row++) {
    $wsData.Range("S$row").Value = $newVersionString
}

$wb.Save()
